# The workbook used to contain three worksheets:
#   Transition_Name_Annot, ISTD_Annot, Sample_Annot
# The new version of MSTemplate_Creator only produces/uses the
# Sample_Annot sheet for this validation test, so drop the other two
# sheets (and let Excel clean up the now-unused shared strings,
# renumber relationship ids, etc. automatically).

$wb = $excel.ActiveWorkbook

# Avoid any "are you sure you want to delete" prompts.
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Transition_Name_Annot").Delete()
$wb.Worksheets.Item("ISTD_Annot").Delete()

# Only Sample_Annot remains - make sure it is the active/selected sheet.
$wb.Worksheets.Item("Sample_Annot").Activate()

$excel.DisplayAlerts = $true
